# Add QUANDL CPI data rows (26-31) to Tabelle1, mirroring the fill-down
# pattern used for the existing "BB" source rows (24-25): KEY/EXPL typed
# per-row, then SOURCE filled for the whole block, then NA_METHOD/LAG/
# HLOC/FREQ copied down from the row above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- KEY (B) / EXPL (C) columns, typed row by row ------------------------
$ws.Range("B26").Value = "RATEINF/CPI_USA"
$ws.Range("C26").Value = "CPI_US"

$ws.Range("B27").Value = "RATEINF/CPI_GBR"
$ws.Range("C27").Value = "CPI_UK"

$ws.Range("B28").Value = "RATEINF/CPI_ITA"
$ws.Range("C28").Value = "CPI_IT"

$ws.Range("B29").Value = "RATEINF/CPI_FRA"
$ws.Range("C29").Value = "CPI_FR"

$ws.Range("B30").Value = "RATEINF/CPI_EUR"
$ws.Range("B31").Value = "RATEINF/CPI_DEU"
$ws.Range("C30").Value = "CPI_EUR"
$ws.Range("C31").Value = "CPI_GER"

# --- SOURCE (A) filled last for the whole block ---------------------------
$ws.Range("A26:A31").Value = "QUANDL"

# --- remaining columns copied down from row 25 (same values as row 25) ---
$ws.Range("D26:D31").Value = "LAST"
$ws.Range("F26:F31").Value = 25
for ($r = 26; $r -le 31; $r++) {
    $ws.Range("G25").Copy($ws.Range("G$r"))
}
$ws.Range("H26:H31").Value = "M"

# --- move the selection to where the user ended up -------------------------
$ws.Range("B27").Select()
